$wb = $excel.ActiveWorkbook

# ===================================================================
# Sheet "Prix Spot": a new day (01-sep) was published by EPEX Spot.
# Add a new column CB mirroring the existing day columns: a header
# in row 1 and 24 hourly prices in rows 2-25.
# ===================================================================
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$colNew = 80   # column CB
$colPrev = 79  # column CA (existing last day column, used as format source)

$wsPrix.Cells.Item(1, $colNew).Value = "01-sep"
# Copy the header's formatting (bold, border, centered) from the
# previous day's header cell so the new column matches the others.
$wsPrix.Cells.Item(1, $colPrev).Copy() | Out-Null
$wsPrix.Cells.Item(1, $colNew).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$prixValues = @(
    35.79,
    18.07,
    22.34,
    17.64,
    15,
    25.21,
    33.23,
    53.18,
    64,
    35.79,
    14.03,
    16.64,
    13.73,
    6.07,
    0,
    4.65,
    8.13,
    14.08,
    23.08,
    55,
    80,
    84.24,
    84.74,
    78.5
)

for ($i = 0; $i -lt $prixValues.Count; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, $colNew).Value = $prixValues[$i]
}

# ===================================================================
# Sheet "Gaz": append the two new daily closing prices.
# ===================================================================
$wsGaz = $wb.Worksheets.Item("Gaz")

function Add-DateRow($ws, $row, $dateText, $price) {
    $cellDate = $ws.Cells.Item($row, 1)
    # Force the date column to remain plain text (matches every other
    # row in the sheet, which stores dates as literal strings, not
    # Excel date serials), then drop back to the default "Normal"
    # style so no stray number-format style gets introduced.
    $cellDate.NumberFormat = "@"
    $cellDate.Value = $dateText
    $cellDate.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $price
}

Add-DateRow $wsGaz 77 "2025-08-30" 30.225
Add-DateRow $wsGaz 78 "2025-08-31" 30.225

# ===================================================================
# Sheet "CO2": append the two new daily closing prices.
# ===================================================================
$wsCo2 = $wb.Worksheets.Item("CO2")

Add-DateRow $wsCo2 77 "2025-08-30" 71.1
Add-DateRow $wsCo2 78 "2025-08-31" 71.1
